$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- "Dialogs" progress table (rows 27-31, column B = Rnd01..Rnd04) ---
# Round 3 dialogs are now finished ("DONE", green), matching the style
# already used by Round 1 / Round 2 (B28/B29). The previous "in progress"
# note is gone from Round 3, and Round 1 / Round 2 no longer show the DONE
# marker (their cells are cleared but keep their green formatting), while
# Round 4 (B31) is cleared back to blank (keeping its formatting).

# Copy the "DONE" look (green fill) from B28 onto B30, then set its text.
$ws.Cells.Item(28, 2).Copy()
$ws.Cells.Item(30, 2).PasteSpecial(-4122)
$ws.Cells.Item(30, 2).Value = "DONE"

# Clear the old DONE markers on Round 1 / Round 2 (formatting stays as-is).
$ws.Cells.Item(28, 2).Value = ""
$ws.Cells.Item(29, 2).Value = ""

# Clear the placeholder text on Round 4 (formatting stays as-is).
$ws.Cells.Item(31, 2).Value = ""

# --- "Undeads" table: the javelin-throw bug note on C25 was resolved, ---
# --- cell goes back to the normal "doplnit" placeholder look (matches D25) ---
$ws.Cells.Item(25, 4).Copy()
$ws.Cells.Item(25, 3).PasteSpecial(-4122)
$ws.Cells.Item(25, 3).Value = "doplnit"

$excel.CutCopyMode = 0

# --- restore the reported selection ---
$ws.Range("B30").Select()
